$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Region being edited: Part 3 - Evaluation, questions 2-4 and their
# answer placeholders (originally paragraphs 50-61, 1-indexed).
#
# We work from the bottom of the region upward so that paragraph
# indices used for the not-yet-processed (earlier) paragraphs stay
# valid while we edit/insert/delete later ones.
# --------------------------------------------------------------------

# --- Step 1: remove the old Q4 block (empty separator, Q4 heading,
#     Q4 answer placeholder) - its content is reinserted lower down,
#     right after the (to become) Q3 text, per the target layout.
$p57 = $d.Paragraphs.Item(57)
$p59 = $d.Paragraphs.Item(59)
$delRange = $d.Range($p57.Range.Start, $p59.Range.End)
$delRange.Delete()

# --- Step 2: paragraph 56 (Heading3 "[1-2 paragraphs, 1 point]" -
#     answer placeholder for what was Q3) - normalize to a single run
#     (content unchanged).
$p56 = $d.Paragraphs.Item(56)
$p56.Range.Text = "[1-2 paragraphs, 1 point]"

# --- Step 3: paragraph 55 (Heading2, was "3)  How could you refine
#     your data analytics?") now becomes question 4.
$p55 = $d.Paragraphs.Item(55)
$p55.Range.Text = "4)  Are there any implications for employers and employees based on the findings you obtained? Justify your answer."

# --- Step 4: paragraph 54 (empty separator between old Q2-answer and
#     old Q3 heading) is replaced by two new body paragraphs discussing
#     refinements to the data analytics.
$p54 = $d.Paragraphs.Item(54)
$p54.Range.Text = "There are a variety of refinements that could be made to improve the quality of the data analysis. Firstly, it would have been beneficial to have used SEEK data from a wider range of time. The dataset provided only included data from the end of 2018 and the start of 2019 and misses half of the months (April, May, June, July, August, September). Since July is not included, it was impossible to see if the Australian end of the financial year impacted the number of job listings. In addition to this, the 2019 data contained Classifications and Location columns with NaN values and thus could not contribute to a large part of the analysis. With a more robust and wider-ranging dataset it would have been possible to analyse the market trends over time. For example, it might have allowed us to investigate the popularity of certain IT technologies over time, illustrating if some are increasing or some are decreasing in popularity (perhaps Java is slowly decreasing, but Python is quickly increasing?)."
$p54.Range.Style = "Normal"

$p54.Range.InsertParagraphAfter()
$p55b = $d.Paragraphs.Item(55)
$p55b.Range.Text = "While analysing the most popular IT technologies, it would have been an improvement to use a separate and more comprehensive dataset containing the names of technologies. In the analysis a self-composed list of 30 technologies was used.  Because it was self-composed and small in size there is a possibility that some technologies have been left out. If this is the case, the data may be misleading as to what technologies should be studied. Additionally, the method of finding these keywords could be improved. There are some cases where technologies may be spelt differently such as “Objective-C” and “Objective C” or “SQL” and “MySQL”. To improve this, more advanced Natural Language Tool Kit (NLTK) processing should be implemented such as stemming."
$p55b.Range.Style = "Normal"
$full = $p55b.Range
$full.Bold = 1
$textOnly = $d.Range($full.Start, $full.End - 1)
$textOnly.Bold = 0
$endRange = $d.Range($full.End - 1, $full.End - 1)
$endRange.InsertBreak(6)

# --- Step 5: paragraph 53 (Heading3 "[1-2 paragraphs, 1 point]" -
#     answer placeholder for old Q2) - normalize to a single run.
$p53 = $d.Paragraphs.Item(53)
$p53.Range.Text = "[1-2 paragraphs, 1 point]"

# --- Step 6: paragraph 52 (Heading2, was "2)  What actions for
#     balancing the markets...") now becomes question 3.
$p52 = $d.Paragraphs.Item(52)
$p52.Range.Text = "3)  How could you refine your data analytics?"

# --- Step 7: paragraph 51 (empty separator) is left untouched.

# --- Step 8: paragraph 50 (Heading3 "[1-2 paragraphs, 2 points]" -
#     answer placeholder for Q1) widens to "2-3 paragraphs".
$p50 = $d.Paragraphs.Item(50)
$p50.Range.Text = "[2-3 paragraphs, 2 points]"

# Insert the new Q2 block right after it: empty separator, the Q2
# heading, its answer placeholder, and a note paragraph.
$p50.Range.InsertParagraphAfter()
$pA = $d.Paragraphs.Item(51)
$pA.Range.Style = "Normal"

$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item(52)
$pB.Range.Text = "2)  What actions for balancing the markets do you suggest based on your findings?"
$pB.Range.Style = "Heading 2"

$pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Item(53)
$pC.Range.Text = "[1-2 paragraphs, 1 point]"
$pC.Range.Style = "Heading 3"

$pC.Range.InsertParagraphAfter()
$pD = $d.Paragraphs.Item(54)
$pD.Range.Text = "? – sent email to Henry asking what this means."
$pD.Range.Style = "Normal"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
